$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update resultado (G) / profit (H) for existing rows ---
$updates = @(
    @{Row=42;  Resultado="Fallo";   Profit=-1},
    @{Row=44;  Resultado="Fallo";   Profit=-1},
    @{Row=64;  Resultado="Acierto"; Profit=2},
    @{Row=73;  Resultado="Acierto"; Profit=1.75},
    @{Row=89;  Resultado="Fallo";   Profit=-1},
    @{Row=90;  Resultado="Acierto"; Profit=1.1},
    @{Row=91;  Resultado="Fallo";   Profit=-1},
    @{Row=93;  Resultado="Fallo";   Profit=-1},
    @{Row=94;  Resultado="Fallo";   Profit=-1},
    @{Row=101; Resultado="Fallo";   Profit=-1},
    @{Row=102; Resultado="Fallo";   Profit=-1},
    @{Row=106; Resultado="Acierto"; Profit=1.62},
    @{Row=111; Resultado="Fallo";   Profit=-1},
    @{Row=112; Resultado="Fallo";   Profit=-1},
    @{Row=113; Resultado="Acierto"; Profit=2.4},
    @{Row=114; Resultado="Fallo";   Profit=-1},
    @{Row=116; Resultado="Fallo";   Profit=-1}
)

foreach ($u in $updates) {
    $r = $u.Row
    $ws.Cells.Item($r, 7).Value = $u.Resultado
    $ws.Cells.Item($r, 8).Value = $u.Profit
}

# --- Append new rows 119 and 120 ---
# Note: the "fecha" column holds text that looks like an ISO date
# ("2025-09-04"); a leading apostrophe forces Excel to keep it as literal
# text instead of auto-converting it to a date serial number, matching the
# source data (stored as a plain string).
$newRows = @(
    @{Row=119; A=14552573; B="'2025-09-04"; C="Carlo Alberto Caniato"; D="Federico Bondioli";  E="Gana Federico Bondioli";   F=2.63},
    @{Row=120; A=14552526; B="'2025-09-04"; C="Giulio Zeppieri";       D="Petr Bar Biryukov";   E="Gana Petr Bar Biryukov";  F=3.75}
)

foreach ($nr in $newRows) {
    $r = $nr.Row
    $ws.Cells.Item($r, 1).Value = $nr.A
    $ws.Cells.Item($r, 2).Value = $nr.B
    $ws.Cells.Item($r, 3).Value = $nr.C
    $ws.Cells.Item($r, 4).Value = $nr.D
    $ws.Cells.Item($r, 5).Value = $nr.E
    $ws.Cells.Item($r, 6).Value = $nr.F
}
